$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet currently has two header rows (row1: units for some cols, row2:
# Hiver/Eté/Année-ish sub headers) followed by 12 data rows (rows 3-14).
# The target layout is a single header row followed by the same 12 data
# rows (now rows 2-13). Deleting row 1 shifts row 2 (the "(m3/s)/(MW)/(GWh)"
# row) up into row 1, and the data rows up into rows 2-13 -- exactly the
# row alignment we need, so we only have to rewrite the header text/styles.
$ws.Rows.Item(1).Delete()

# New leading identifier/name/date columns (no custom style -> default s=0).
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

# F1 keeps the existing "(m3/s)" shared string/style; the remaining unit
# headers are replaced with the new, more descriptive labels.
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# F1:K1 all carry a dedicated header style: same font as the data cells
# (Arial 9, cellXf fontId=1) but General number format with no explicit
# applyNumberFormat flag. Create a transient named style to mint that cellXf,
# apply it to the header row, then remove the named style again so the
# workbook's cellStyles/cellStyleXfs tables stay as they were -- only a new
# cellXf entry remains behind, referenced by the header cells.
$headerStyle = $wb.Styles.Add("TmpHeaderStyle")
$headerStyle.Font.Name = "Arial"
$headerStyle.Font.Size = 9
$ws.Range("F1:K1").Style = "TmpHeaderStyle"
$wb.Styles.Item("TmpHeaderStyle").Delete()

# Restore the original selection/active-cell marker for the new header row.
$ws.Range("A2:K2").Select() | Out-Null
